# edit.ps1
# Applies "edit route slug di Admin" change:
#   - Row 3 (Kaos Panjang / kaos-panjang) -> (Lengen Panjang / lengen-panjang)
#   - Row 6 (Kaos Kaki / kaos kaki)       -> (Aksesoris / aksesoris)
#   - Update the active selection on the sheet to H3 (and reset scroll to top)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the category name + slug for the two renamed rows
$ws.Range("B3").Value = "Lengen Panjang"
$ws.Range("C3").Value = "lengen-panjang"

$ws.Range("B6").Value = "Aksesoris"
$ws.Range("C6").Value = "aksesoris"

# Reset scroll position to the top-left and move/select H3
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H3").Select()
